$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q)
$ws.Columns("N").Insert()

# This sheet becomes the active sheet/tab, with a new selection
$ws.Activate()
$ws.Range("S6").Select()

Write-Output "done"
